# "some results week over week" - fill in this week's answer_truth results
# for the Game of Thrones answer_structure sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("answer_structure")

# --- 1. Characters who die (switch from the default "Lives" to "Dies") ---
# 08. Cersei Lannister, 10. Jaime Lannister, 11. Jorah Mormont,
# 13. Theon Greyjoy, 15. The Hound, 18. Melisandre, 19. Missandei, 22. Varys
$ws.Range("E9").Value  = "Dies"
$ws.Range("E11").Value = "Dies"
$ws.Range("E12").Value = "Dies"
$ws.Range("E14").Value = "Dies"
$ws.Range("E16").Value = "Dies"
$ws.Range("E19").Value = "Dies"
$ws.Range("E20").Value = "Dies"
$ws.Range("E23").Value = "Dies"

# --- 2. Newly-graded questions: mark as included and fill the answer_truth ---
# 29. Pick a character that kills a White Walker -> Arya Stark
$ws.Range("C30").Value = $true
$ws.Range("E30").Value = "Arya Stark"

# 31/32/33/35/36. Does Arya personally kill ... -> No
$ws.Range("C32").Value = $true
$ws.Range("E32").Value = "No"
$ws.Range("C33").Value = $true
$ws.Range("E33").Value = "No"
$ws.Range("C34").Value = $true
$ws.Range("E34").Value = "No"
$ws.Range("C36").Value = $true
$ws.Range("E36").Value = "No"
$ws.Range("C37").Value = $true
$ws.Range("E37").Value = "No"

# 30. Pick a character that becomes reanimated as a wight -> Edd Tollett
$ws.Range("C31").Value = $true
$ws.Range("E31").Value = "Edd Tollett"

# 41. Pick a character that returns after their brief hiatus -> Rhaegal
$ws.Range("C45").Value = $true
$ws.Range("E45").Value = "Rhaegal"

# 26. Which major character kills which major character? -> Arya kills the Night King
$ws.Range("C27").Value = $true
$ws.Range("E27").Value = "Arya kills the Knight's King"

# 27. Which TWO supporting characters kill which TWO supporting characters?
$ws.Range("C28").Value = $true
$ws.Range("E28").Value = "The Hound kills Ser Gregor, Ser Gregor kills Qyburn"

# 42. Will the following characters appear naked -> No
$ws.Range("C46").Value = $true
$ws.Range("E46").Value = "No"

# --- 3. multiple_answers flips to No for the dragon-flying questions ---
$ws.Range("D38").Value = $false
$ws.Range("D39").Value = $false
$ws.Range("D40").Value = $false
$ws.Range("D41").Value = $false

# --- 4. Question 49 gains a multiple_answers=No flag ---
$ws.Range("D50").Value = $false

# --- 5. Sort the question table by question (A2:E50); data is already in
#        order so this only stamps the "remembered sort" state ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A50"))
$ws.Sort.SetRange($ws.Range("A2:E50"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- 6. Leave the selection like the author did before saving ---
$ws.Range("C2:E52").Select()
